$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 2 entirely (the "com.singleton.strechy / taxi game /
# galiatia942@gmail.com / syechimovitz@gmail.com ..." review), shifting
# all subsequent rows up by one.
$ws.Rows.Item(2).Delete()

# Row deletion doesn't re-anchor the worksheet's mailto hyperlinks, so
# rebuild them pointing one row higher than before, preserving order
# (and thus the relationship-id sequence) and display text.
$ws.Hyperlinks.Delete()

# Cell C4 keeps the plain "data" formatting (style index 2) and never
# itself becomes a hyperlink target below, so use it as a formatting
# donor: Hyperlinks.Add() stamps its own (underlined / themed) font on
# the target cell, and we paste the donor's formatting back over it so
# the cell keeps looking like ordinary data, same as the other data
# cells in this column.
$fmtDonor = $ws.Range("C4")

function Add-PlainHyperlink($addr, $mailto, $display) {
    $target = $ws.Range($addr)
    $ws.Hyperlinks.Add($target, $mailto, "", "", $display)
    $fmtDonor.Copy()
    $target.PasteSpecial(-4122) # xlPasteFormats
}

Add-PlainHyperlink "D2" "mailto:shmulmaor2@gmail.com" "shmulmaor2@gmail.com"
Add-PlainHyperlink "C3" "mailto:rocketaso@gmail.com" "rocketaso@gmail.com"
Add-PlainHyperlink "D3" "mailto:armonravid@gmail.com" "armonravid@gmail.com"
Add-PlainHyperlink "C5" "mailto:ronoren61@gmail.com" "ronoren61@gmail.com"
Add-PlainHyperlink "D5" "mailto:nitanoren23@gmail.com" "nitanoren23@gmail.com"
Add-PlainHyperlink "C7" "mailto:danfogel100@gmail.com" "danfogel100@gmail.com"
Add-PlainHyperlink "D7" "mailto:avishaybar12@gmail.com" "avishaybar12@gmail.com"
Add-PlainHyperlink "C8" "mailto:danfogel100@gmail.com" "danfogel100@gmail.com"
Add-PlainHyperlink "D8" "mailto:avishaybar12@gmail.com" "avishaybar12@gmail.com"
Add-PlainHyperlink "D9" "mailto:jorjkluni03@gmail.com" "jorjkluni03@gmail.com"

$excel.CutCopyMode = $false
